# Update the "two-digit number divided by one-digit number" practice
# table. Each populated table row holds five division problems; replace
# the problem text in each cell with the new problem text, cell by cell,
# so that duplicate source values (e.g. "69÷2=" appearing twice) are each
# replaced with their own distinct target value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map: table row (1-based) -> array of 5 new values for columns 1..5
$updates = @{
    1  = @("56÷7=", "49÷8=", "12÷5=", "69÷4=", "15÷3=")
    5  = @("68÷7=", "37÷3=", "91÷5=", "12÷5=", "63÷5=")
    9  = @("20÷2=", "89÷6=", "92÷3=", "90÷3=", "41÷4=")
    13 = @("25÷2=", "30÷2=", "37÷5=", "74÷5=", "43÷4=")
    17 = @("67÷2=", "88÷3=", "50÷2=", "93÷2=", "35÷3=")
}

foreach ($rowIndex in $updates.Keys) {
    $values = $updates[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
